# Added Time Logs testcases
# a. Created 2 Time Logs TestCases
# b. edited Homepage
# c. Added EditTimeLogs

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Logs")
$ws.Activate()

# Update row 3 of the "Time Logs" sheet with the second Time Logs testcase
$ws.Range("A3").Value = "TC002_TimeTracker_TimeLogs_LogTimeOut"
$ws.Range("B3").Value = "glenn.mamaril"

# Update the selection to C3 to match the saved view state
$ws.Range("C3").Select()
